$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Parametros fractales" sheet: the strategy table (header in row 1, kept
# as-is) is extended from 3 data rows (2 strategies + FIRX summary) to 9
# strategies + the FIRX summary row, which moves down from row 4 to row 11.
# Two of the new strategies (rows 9-10) introduce a new "OPEN" reference
# value that did not exist in the workbook before.
# Columns: A..M
$data = @(
    @(1, "CLOSE", ">", "max", "HIGH", "BUY", "CLOSE", "<", "min", "LOW", "SELL", -4, -1),
    @(2, "CLOSE", ">", "max", "HIGH", "BUY", "CLOSE", "<", "min", "LOW", "SELL", -5, -2),
    @(3, "CLOSE", ">", "min", "LOW", "BUY", "CLOSE", "<", "max", "HIGH", "SELL", -4, -1),
    @(4, "CLOSE", ">", "min", "LOW", "BUY", "CLOSE", "<", "max", "HIGH", "SELL", -5, -1),
    @(5, "CLOSE", ">", "min", "LOW", "BUY", "CLOSE", "<", "max", "HIGH", "SELL", -6, -1),
    @(6, "CLOSE", ">", "min", "LOW", "BUY", "CLOSE", "<", "max", "HIGH", "SELL", -3, -1),
    @(7, "CLOSE", ">", "min", "LOW", "BUY", "CLOSE", "<", "max", "HIGH", "SELL", -6, -2),
    @(8, "CLOSE", ">", "max", "OPEN", "BUY", "CLOSE", "<", "min", "LOW", "SELL", -4, -1),
    @(9, "CLOSE", ">", "max", "OPEN", "BUY", "CLOSE", "<", "min", "HIGH", "SELL", -4, -1)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $ws.Cells.Item($row, 9).Value = $r[8]
    $ws.Cells.Item($row, 10).Value = $r[9]
    $ws.Cells.Item($row, 11).Value = $r[10]
    $ws.Cells.Item($row, 12).Value = $r[11]
    $ws.Cells.Item($row, 13).Value = $r[12]
    $row++
}

# Move the "FIRX" summary row from row 4 to row 11, keep same content
$ws.Cells.Item(11, 1).Value = "FIRX"
$ws.Cells.Item(11, 2).Value = "FIX"
$ws.Cells.Item(11, 3).Value = ">"
$ws.Cells.Item(11, 4).Value = "FDX"
$ws.Cells.Item(11, 5).Value = "FDX"
$ws.Cells.Item(11, 6).Value = "BUY"
$ws.Cells.Item(11, 7).Value = "FIX"
$ws.Cells.Item(11, 8).Value = "<"
$ws.Cells.Item(11, 9).Value = "FDX"
$ws.Cells.Item(11, 10).Value = "FDX"
$ws.Cells.Item(11, 11).Value = "SELL"

# Apply centered style (style index 1 in target) to the new rows 5-11
$ws.Range("A5:M11").HorizontalAlignment = -4108  # xlCenter

# Match the workbook's recorded selection (A2) after editing
$ws.Range("A2").Select() | Out-Null
